# Updated main GSC export data: the oldest day (2025-11-02) row is removed
# from the "Chart" sheet, causing every subsequent day's row to shift up by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Delete row 2 (the 2025-11-02 data row); rows below shift up automatically.
$ws.Rows("2:2").Delete()
